$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D2" = 3261
    "E2" = 258
    "F2" = 258
    "G2" = 298
    "H2" = 225
    "I2" = 220
    "J2" = 5
    "K2" = 3596
    "L2" = 1002
    "M2" = 2594
    "N2" = 2493
    "O2" = 101
    "P2" = 90
    "Q2" = 434
    "R2" = -422
    "S2" = -31
    "T2" = 203
    "U2" = 231
    "V2" = 383
    "W2" = 7.9
    "X2" = 6.9
    "Y2" = 9.15
    "Z2" = 6.28
    "AA2" = 38.65
    "AB2" = 2799.31
    "AC2" = 1226
    "AD2" = 10.07
    "AE2" = 15135
    "AF2" = 0.82
    "AH2" = 3.24
    "AI2" = 30.01
    "AJ2" = 17900000
    "D3" = 2673
    "E3" = 242
    "F3" = 242
    "G3" = 265
    "H3" = 208
    "I3" = 199
    "J3" = 9
    "K3" = 3993
    "L3" = 1244
    "M3" = 2749
    "N3" = 2653
    "O3" = 96
    "P3" = 90
    "Q3" = 495
    "R3" = -361
    "S3" = -96
    "T3" = 300
    "U3" = 196
    "V3" = 368
    "W3" = 9.050000000000001
    "X3" = 7.78
    "Y3" = 7.72
    "Z3" = 5.48
    "AA3" = 45.24
    "AB3" = 2931.64
    "AC3" = 1110
    "AD3" = 13.74
    "AE3" = 16108
    "AF3" = 0.95
    "AG3" = 400
    "AH3" = 2.62
    "AI3" = 33.16
    "AJ3" = 17900000
    "D4" = 2985
    "E4" = 193
    "F4" = 193
    "G4" = 197
    "H4" = 149
    "I4" = 145
    "J4" = 5
    "K4" = 4268
    "L4" = 1421
    "M4" = 2847
    "N4" = 2747
    "O4" = 100
    "P4" = 90
    "Q4" = 504
    "R4" = -568
    "S4" = 54
    "T4" = 617
    "U4" = -114
    "V4" = 496
    "W4" = 6.48
    "X4" = 5
    "Y4" = 5.36
    "Z4" = 3.61
    "AA4" = 49.92
    "AB4" = 3043.5
    "AC4" = 808
    "AD4" = 21.1
    "AE4" = 16325
    "AF4" = 1.04
    "AG4" = 400
    "AH4" = 2.35
    "AI4" = 46.55
    "AJ4" = 17900000
    "D5" = 4456
    "E5" = 279
    "F5" = 279
    "G5" = 300
    "H5" = 249
    "I5" = 237
    "J5" = 12
    "K5" = 4608
    "L5" = 1772
    "M5" = 2836
    "N5" = 2836
    "O5" = 0
    "P5" = 90
    "Q5" = 384
    "R5" = -790
    "S5" = 552
    "T5" = 768
    "U5" = -384
    "V5" = 1138
    "W5" = 6.26
    "X5" = 5.59
    "Y5" = 8.48
    "Z5" = 5.61
    "AA5" = 62.46
    "AB5" = 3190.93
    "AC5" = 1322
    "AD5" = 20.16
    "AE5" = 16853
    "AF5" = 1.58
    "AG5" = 450
    "AH5" = 1.69
    "AI5" = 32
    "AJ5" = 17900000
    "D6" = 3583
    "E6" = -77
    "F6" = -77
    "G6" = -14
    "H6" = -13
    "I6" = -13
    "K6" = 3522
    "L6" = 752
    "M6" = 2770
    "N6" = 2770
    "P6" = 90
    "Q6" = 365
    "R6" = 315
    "S6" = -943
    "T6" = 298
    "U6" = 67
    "V6" = 313
    "W6" = -2.15
    "X6" = -0.36
    "Y6" = -0.46
    "Z6" = -0.32
    "AA6" = 27.17
    "AB6" = 3123.26
    "AC6" = -72
    "AD6" = -164.61
    "AE6" = 16457
    "AF6" = 0.72
    "AI6" = -262.32
    "AJ6" = 17900000
    "D7" = 3880
    "E7" = -26
    "G7" = 130
    "H7" = 75
    "I7" = 75
    "K7" = 3550
    "L7" = 740
    "M7" = 2810
    "N7" = 2810
    "P7" = 90
    "Q7" = 290
    "R7" = -210
    "S7" = -90
    "U7" = 180
    "W7" = -0.67
    "X7" = 1.93
    "Y7" = 2.69
    "Z7" = 2.12
    "AA7" = 26.33
    "AC7" = 419
    "AD7" = 23.87
    "AE7" = 16697
    "AF7" = 0.6
    "AG7" = 200
    "AH7" = 2
    "AI7" = 47.73
    "D8" = 4563
    "E8" = 185
    "G8" = 308
    "H8" = 241
    "I8" = 241
    "K8" = 3730
    "L8" = 710
    "M8" = 3020
    "N8" = 3020
    "P8" = 90
    "Q8" = 460
    "R8" = -190
    "S8" = -120
    "U8" = 330
    "W8" = 4.05
    "X8" = 5.28
    "Y8" = 8.27
    "Z8" = 6.62
    "AA8" = 23.51
    "AC8" = 1346
    "AD8" = 7.43
    "AE8" = 17945
    "AF8" = 0.5600000000000001
    "AG8" = 200
    "AH8" = 2
    "AI8" = 14.85
    "D9" = 4840
    "E9" = 230
    "G9" = 360
    "H9" = 280
    "I9" = 280
    "K9" = 3970
    "L9" = 710
    "M9" = 3260
    "N9" = 3260
    "P9" = 90
    "Q9" = 580
    "R9" = -470
    "S9" = -70
    "U9" = 160
    "W9" = 4.75
    "X9" = 5.79
    "Y9" = 8.92
    "Z9" = 7.27
    "AA9" = 21.78
    "AC9" = 1564
    "AD9" = 6.39
    "AE9" = 19371
    "AF9" = 0.52
    "AG9" = 200
    "AH9" = 2
    "AI9" = 12.79
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# Cells removed entirely in the target (matches <c> tag deletion in OOXML)
$cellsToClear = @(
    "AG6",
    "AH6",
    "T7",
    "T8",
    "T9",
)
foreach ($ref in $cellsToClear) {
    $ws.Range($ref).ClearContents()
}
